$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.748.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.805.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5923"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06837"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07495"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.809.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.777"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6248"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.050.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009306"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "75.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.681.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.483"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "211.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.851"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.43%  "
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.893"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1271"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.423"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06193"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.425"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.789"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.766"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.730"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.067"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6429"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.493"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.724"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.560"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01715"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.145.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8797"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.48%  "
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.960.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.01%  "
$ws.Range("E47").Value = "  -1.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.600"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.401"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05468"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4486"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.51%  "
